$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: item id "2940803" stored as text, left-aligned + bordered style ---
$ws.Range("A34").Borders.LineStyle = 1
$ws.Range("A34").HorizontalAlignment = -4131
$ws.Range("A34").Formula = '="2940803"'
$ws.Range("A34").Copy()
$ws.Range("A34").PasteSpecial(-4163)   # xlPasteValues - collapse formula to literal text value

$ws.Range("B34").Value = "LE CHAT ROSE 1L"
$ws.Range("C34").Value = 10
$ws.Range("D34").Value = 368
$ws.Range("B34:D34").Borders.LineStyle = 1

# --- Row 35: item id 2940804 stored as a number, same bordered/left-aligned style ---
$ws.Range("A35").Value = 2940804
$ws.Range("A35").Borders.LineStyle = 1
$ws.Range("A35").HorizontalAlignment = -4131

$ws.Range("B35").Value = "LE CHAT power gel 4L"
$ws.Range("C35").Value = 4
$ws.Range("D35").Value = 1170
$ws.Range("B35:D35").Borders.LineStyle = 1

# Update selection to match the new extent of data
$ws.Range("D36").Select()
